# Apply the edit described by the diff:
# A new weekly price record is inserted as row 604 (pushing the existing
# rows 604-680 down to 605-681), then the new row 604 is populated with
# its own data (date 2023-10-13 / serial 45212, Primera quality, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 604, shifting existing data (604:680) down to (605:681)
$ws.Rows("604:604").Insert()

# Populate the newly inserted row 604 with the new record's values
$ws.Range("A604").Value = 9
$ws.Range("B604").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C604").Value = "Metropolitana"
$ws.Range("D604").Value = 45212
$ws.Range("E604").Value = 13
$ws.Range("F604").Value = 100112012
$ws.Range("G604").Value = "Espinaca"
$ws.Range("H604").Value = "Sin especificar"
$ws.Range("I604").Value = "Primera"
$ws.Range("J604").Value = 160
$ws.Range("K604").Value = 8000
$ws.Range("L604").Value = 9000
$ws.Range("M604").Value = 8500
$ws.Range("N604").Value = '$/cuna 10 kilos'
$ws.Range("O604").Value = "Provincia de Chacabuco"
$ws.Range("P604").Value = 850
$ws.Range("Q604").Value = 10
$ws.Range("R604").Value = "Hortaliza"

# Match the date cell's number format (yyyy-mm-dd hh:mm:ss) used by the rest of column D
$ws.Range("D604").NumberFormat = $ws.Range("D605").NumberFormat
